$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.854.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.51%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.211.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.37%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'604.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.93%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'153.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.211.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.37%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.08%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.65%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.507"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.87%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'38.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.07%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.740.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'Polkadot"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'7.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.21%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'66.074.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.67%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.240.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.23%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.05%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'509.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'15.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.46%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'15.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.22%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'8.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.25%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'85.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.60%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.45%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +3.84%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +3.72%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'28.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.67%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.12%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'55.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.20%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0902"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'476.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.05%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0418"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.15%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.47%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'8.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.296"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.34%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.48%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.950.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.62%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +1.29%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0₃0637"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.30%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'28.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.03%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.05%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +2.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'121.13"
$ws.Range("D51").Style = "Normal"
